$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.417.21"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "3.314.89"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.88"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.31"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  +4.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.61"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "3.833.25"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.49"
$ws.Range("E14").Value = "  +4.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.25"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "3.322.11"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.04"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "59.182.19"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.03"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "304.70"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.25"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.47"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  +6.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.82"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.38"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  +3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.57"
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.69"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  +5.53%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.13"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.38"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.73"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.93"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.81"
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.26"
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.27"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "2.203.00"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.39"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -4.50%  "
